# Update the description text (column B) for several rows on Sheet1
# to reflect the revised test-case wording from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value  = "Input is not accepted. Course name below minimum of 2 characters on line 7."
$ws.Range("B6").Value  = "Input is not accepted.  Course name below minimum of 2 characters on line 2."
$ws.Range("B7").Value  = "Input is not accepted. Course name exceeds maximum of 4 characters on line 11."
$ws.Range("B8").Value  = "Input is not accepted.  Course name exceeds maximum of 4 characters on line 6."
$ws.Range("B27").Value = "Input is not accepted.  Multiple errors encountered. 0 room size in line 2.  Room capacity exceeds maximum in line 5.  Missing section count in line 8. Hours exceeds maximum in line 13."
$ws.Range("B34").Value = "Input is not accepted.  Course number is incorrect on line 1."

# Update the view so the window is scrolled near the bottom of the list
# and the last-edited cell (B34) is the active selection.
$ws.Activate()
$ws.Range("B34").Select()
$excel.ActiveWindow.ScrollRow = 39
